$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Cell VALUES first, in the exact order the new text was authored, so the
#    shared-string table grows with the same index assignment as the target.
# ---------------------------------------------------------------------------

# Second "Clonar un repositorio..." block header/table (rows 27-33)
$ws.Range("A27").Value = "Paso"
$ws.Range("B27").Value = "Comando"
$ws.Range("C27").Value = "Descripción"

$ws.Range("A28").Value = 1
$ws.Range("B28").Value = "git clone https://github.com/YOUR-USERNAME/YOUR-REPOSITORY.git"
$ws.Range("C28").Value = "Clona el repositorio remoto en tu PC"

$ws.Range("A29").Value = 2
$ws.Range("B29").Value = "cd YOUR-REPOSITORY"
$ws.Range("C29").Value = "Cambia al directorio del repositorio clonado"

$ws.Range("A30").Value = 3
$ws.Range("B30").Value = "git status"
$ws.Range("C30").Value = "Muestra el estado de los archivos en el repositorio local"

$ws.Range("A31").Value = 4
$ws.Range("B31").Value = "git add ."
$ws.Range("C31").Value = "Agrega todos los archivos modificados al área de preparación"

$ws.Range("A32").Value = 5
$ws.Range("B32").Value = "git commit -m ""mensaje"""
$ws.Range("C32").Value = "Confirma los cambios con un mensaje"

$ws.Range("A33").Value = 6
$ws.Range("B33").Value = "git push origin BRANCH-NAME"
$ws.Range("C33").Value = "Sube los cambios al repositorio remoto"

# Third block header/table (rows 37-41)
$ws.Range("A37").Value = "Paso"
$ws.Range("B37").Value = "Comando"
$ws.Range("C37").Value = "Descripción"

$ws.Range("A38").Value = 1
$ws.Range("B38").Value = "pwd"
$ws.Range("C38").Value = "Obtiene la ruta del repositorio original"

$ws.Range("A39").Value = 2
$ws.Range("B39").Value = "mkdir proyecto-clonado"
$ws.Range("C39").Value = "Crea un nuevo directorio llamado proyecto-clonado"

$ws.Range("A40").Value = 3
$ws.Range("B40").Value = "cd proyecto-clonado"
$ws.Range("C40").Value = "Cambia al nuevo directorio"

$ws.Range("A41").Value = 4
$ws.Range("B41").Value = "git clone /home/user/proyecto-original"
$ws.Range("C41").Value = "Clona el repositorio original en el nuevo directorio"

# Section titles (new shared string is minted last, matching the source file)
$ws.Range("A26").Value = "Clonar un repositorio de GitHub hacia la PC"
$ws.Range("A36").Value = "Clonar un repositorio de GitHub hacia la PC"

# ---------------------------------------------------------------------------
# 2) Formatting: clone the look of the existing tables by copying formats
#    from equivalent existing cells (keeps the same style entries instead of
#    minting near-duplicate ones).
# ---------------------------------------------------------------------------

# Section-title rows (A/B merged banner), based on A8:B8
$ws.Range("A8").Copy()
$ws.Range("A26").PasteSpecial(-4122)
$ws.Range("B8").Copy()
$ws.Range("B26").PasteSpecial(-4122)
$ws.Range("A8").Copy()
$ws.Range("A36").PasteSpecial(-4122)
$ws.Range("B8").Copy()
$ws.Range("B36").PasteSpecial(-4122)

# Header rows (Paso / Comando / Descripción), based on row 9
$ws.Range("A9").Copy()
$ws.Range("A27").PasteSpecial(-4122)
$ws.Range("B9").Copy()
$ws.Range("B27").PasteSpecial(-4122)
$ws.Range("C9").Copy()
$ws.Range("C27").PasteSpecial(-4122)

$ws.Range("A9").Copy()
$ws.Range("A37").PasteSpecial(-4122)
$ws.Range("B9").Copy()
$ws.Range("B37").PasteSpecial(-4122)
$ws.Range("C9").Copy()
$ws.Range("C37").PasteSpecial(-4122)

# Data rows block 1 (28-33): column A uses the numbered-step style, column B
# alternates between the two data styles, column C is uniform.
$ws.Range("A10").Copy()
$ws.Range("A28:A33").PasteSpecial(-4122)

$ws.Range("B10").Copy()
$ws.Range("B28").PasteSpecial(-4122)
$ws.Range("B11").Copy()
$ws.Range("B29").PasteSpecial(-4122)
$ws.Range("B11").Copy()
$ws.Range("B30").PasteSpecial(-4122)
$ws.Range("B10").Copy()
$ws.Range("B31").PasteSpecial(-4122)
$ws.Range("B11").Copy()
$ws.Range("B32").PasteSpecial(-4122)
$ws.Range("B11").Copy()
$ws.Range("B33").PasteSpecial(-4122)

$ws.Range("C10").Copy()
$ws.Range("C28:C33").PasteSpecial(-4122)

# Data rows block 2 (38-41)
$ws.Range("A10").Copy()
$ws.Range("A38:A41").PasteSpecial(-4122)

$ws.Range("B10").Copy()
$ws.Range("B38").PasteSpecial(-4122)
$ws.Range("B11").Copy()
$ws.Range("B39").PasteSpecial(-4122)
$ws.Range("B11").Copy()
$ws.Range("B40").PasteSpecial(-4122)
$ws.Range("B10").Copy()
$ws.Range("B41").PasteSpecial(-4122)

$ws.Range("C10").Copy()
$ws.Range("C38:C41").PasteSpecial(-4122)

# Trailing blank row 42 mirrors the numbered-row styling with empty content
$ws.Range("A10").Copy()
$ws.Range("A42").PasteSpecial(-4122)
$ws.Range("B11").Copy()
$ws.Range("B42").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C42").PasteSpecial(-4122)

# Two blank spacer rows (34-35) between the two tables: copy the bordered
# data-row formats, then strip the border, mirroring the three new styles
# added to cellXfs.
$ws.Range("A10").Copy()
$ws.Range("A34").PasteSpecial(-4122)
$ws.Range("A34").Borders.LineStyle = -4142

$ws.Range("B11").Copy()
$ws.Range("B34").PasteSpecial(-4122)
$ws.Range("B34").Borders.LineStyle = -4142

$ws.Range("C10").Copy()
$ws.Range("C34").PasteSpecial(-4122)
$ws.Range("C34").Borders.LineStyle = -4142

$ws.Range("A10").Copy()
$ws.Range("A35").PasteSpecial(-4122)
$ws.Range("A35").Borders.LineStyle = -4142

$ws.Range("C10").Copy()
$ws.Range("B35").PasteSpecial(-4122)
$ws.Range("B35").Borders.LineStyle = -4142

$ws.Range("C10").Copy()
$ws.Range("C35").PasteSpecial(-4122)
$ws.Range("C35").Borders.LineStyle = -4142

# ---------------------------------------------------------------------------
# 3) Merge the two new section-title rows
# ---------------------------------------------------------------------------
$ws.Range("A26:B26").Merge()
$ws.Range("A36:B36").Merge()

# ---------------------------------------------------------------------------
# 4) Page setup + view state
# ---------------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("B43").Select()
